$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" (column F) values for the rows noted in the repull/push of data.
$ws.Range("F2").Value = 5
$ws.Range("F3").Value = -10
$ws.Range("F4").Value = 5
$ws.Range("F7").Value = -10
$ws.Range("F10").Value = -4
$ws.Range("F17").Value = -4
$ws.Range("F19").Value = -13
$ws.Range("F21").Value = -6
$ws.Range("F26").Value = -3
$ws.Range("F28").Value = -5
$ws.Range("F30").Value = -3
$ws.Range("F33").Value = 0
